# Bug Fix: Contact form and OTP generation
# Adds two new contact-form submissions (rows 11 and 12) to the Contacts sheet,
# matching the entries created while fixing the contact form / OTP flow.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contacts")

# Row 11: Richard, resubmission with blank message (OTP flow)
$ws.Range("A11").Value = "Richard "
$ws.Range("B11").Value = "9854747474"
$ws.Range("C11").Value = "2BHK"
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = "11/11/2025, 10:26:01 pm"

# Row 12: Robert, new submission with blank message (OTP flow)
$ws.Range("A12").Value = "Robert"
$ws.Range("B12").Value = "9854747474"
$ws.Range("C12").Value = "2BHK"
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = "11/11/2025, 10:54:36 pm"
